$d = $word.ActiveDocument

# Helper: insert `$text` right before the trailing paragraph-mark of
# paragraph `$p`, as a brand-new bold run (kept distinct from any
# neighbouring run — even one with identical formatting — by forcing a
# formatting toggle immediately after the insert).
function Add-BoldRun($p, [string]$text) {
    $r = $p.Range
    $insertPoint = $d.Range($r.End - 1, $r.End - 1)
    $insertPoint.InsertAfter($text)
    # Toggling Bold off then back on forces the engine to keep this
    # text in its own <w:r> instead of folding it back into the
    # preceding run that already carries <w:b/>.
    $insertPoint.Font.Bold = 0
    $insertPoint.Font.Bold = 1
}

$nameHits = 0
$idHits = 0

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    if ($t -match "^Name:\r?$") {
        $nameHits = $nameHits + 1
        if ($nameHits -eq 3) {
            Add-BoldRun $p " Carol Rameder"
        }
    }
    elseif ($t -match "^Student id:\r?$") {
        $idHits = $idHits + 1
        if ($idHits -eq 3) {
            Add-BoldRun $p " "
            Add-BoldRun $p "crr940"
        }
    }
}
